# Week 2 PPT update
#
# 1. Slide 1 (Title shape): remove the stray empty run that precedes the
#    first manual line break in the title text box.
# 2. Slide 2 (Content Placeholder): append a sentence to the "SolidWorks
#    API" bullet.
# 3. Slide 6 (Content Placeholder): change "README" to "guide" in the
#    "Document everything thoroughly" bullet.

$p = $ppt.ActivePresentation

# --- 1. Slide 1: drop the leading empty run in the title -------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item("Title 1")
$titleRange = $titleShape.TextFrame.TextRange

# The leading run has no characters in it, so it can't be targeted with a
# Characters(start, length) selection directly. Instead, stamp a temporary
# marker character into it (InsertBefore lands inside that empty run), then
# delete just that marker character - which removes the now-empty run node
# entirely instead of leaving a hollow <a:r><a:t/></a:r> behind.
[void]$titleRange.InsertBefore("@@TMP@@")
$markerRange = $titleRange.Characters(1, 7)
$markerRange.Delete()

# --- 2. Slide 2: extend the SolidWorks API bullet ---------------------------
$slide2 = $p.Slides.Item(2)
$bodyShape2 = $slide2.Shapes.Item("Content Placeholder 2")
$bodyRange2 = $bodyShape2.TextFrame.TextRange

$oldSnippet = " (custom API) and the official SolidWorks API"
$newSnippet = " (custom API) and the official SolidWorks API. More to be done here!"
$fullText2 = $bodyRange2.Text
$pos2 = $fullText2.IndexOf($oldSnippet)
$target2 = $bodyRange2.Characters($pos2 + 1, $oldSnippet.Length)
$target2.Text = $newSnippet

# --- 3. Slide 6: README -> guide --------------------------------------------
$slide6 = $p.Slides.Item(6)
$bodyShape6 = $slide6.Shapes.Item("Content Placeholder 2")
$bodyRange6 = $bodyShape6.TextFrame.TextRange

$oldLine = "Document everything thoroughly. This project is confusing to set-up, so I will make a README"
$newLine = "Document everything thoroughly. This project is confusing to set-up, so I will make a guide"
$fullText6 = $bodyRange6.Text
$pos6 = $fullText6.IndexOf($oldLine)
$target6 = $bodyRange6.Characters($pos6 + 1, $oldLine.Length)
$target6.Text = $newLine
